$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: update the Handoff/Handback datetimes for the second data row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-16 14:58:57"
$wsZhCn.Range("G3").Value = "2016-02-16 14:59:59"

# "de-de" worksheet: update the Handoff/Handback datetimes for the second data row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-16 14:59:11"
$wsDeDe.Range("G3").Value = "2016-02-16 15:00:34"
